$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the whole editable range to Text format first so that numeric-looking
# strings (e.g. "0.9991", "1.847.60") are stored as text, matching the source
# inlineStr cells instead of being auto-converted to numbers by Excel.
$editRange = $ws.Range("B2:E51")
$editRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "30.202.37"
$ws.Range("E2").Value = "  -0.69%  "

# Row 3
$ws.Range("D3").Value = "1.851.34"
$ws.Range("E3").Value = "  -2.05%  "

# Row 4
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
$ws.Range("D5").Value = "236.36"
$ws.Range("E5").Value = "  -0.64%  "

# Row 6
$ws.Range("D6").Value = "0.9992"
$ws.Range("E6").Value = "  -0.19%  "

# Row 7
$ws.Range("D7").Value = "0.4779"
$ws.Range("E7").Value = "  -2.32%  "

# Row 8
$ws.Range("D8").Value = "0.2808"
$ws.Range("E8").Value = "  -4.14%  "

# Row 9
$ws.Range("D9").Value = "0.06478"
$ws.Range("E9").Value = "  -3.18%  "

# Row 10
$ws.Range("D10").Value = "1.857.88"
$ws.Range("E10").Value = "  -1.55%  "

# Row 11
$ws.Range("D11").Value = "0.07314"
$ws.Range("E11").Value = "  -0.44%  "

# Row 12
$ws.Range("D12").Value = "16.29"
$ws.Range("E12").Value = "  -3.88%  "

# Row 13
$ws.Range("D13").Value = "5.112"
$ws.Range("E13").Value = "  -0.22%  "

# Row 14
$ws.Range("D14").Value = "87.20"
$ws.Range("E14").Value = "  -0.45%  "

# Row 15
$ws.Range("D15").Value = "0.6467"
$ws.Range("E15").Value = "  -2.36%  "

# Row 16
$ws.Range("D16").Value = "30.141.14"
$ws.Range("E16").Value = "  -0.83%  "

# Row 17
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "13.24"
$ws.Range("E17").Value = "  -1.39%  "

# Row 18
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "0.9995"
$ws.Range("E18").Value = "  -0.19%  "

# Row 19
$ws.Range("D19").Value = "0.000007633"
$ws.Range("E19").Value = "  -2.39%  "

# Row 20
$ws.Range("D20").Value = "225.63"
$ws.Range("E20").Value = "  +18.67%  "

# Row 21
$ws.Range("D21").Value = "2.097.53"
$ws.Range("E21").Value = "  -3.13%  "

# Row 22
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "5.293"
$ws.Range("E22").Value = "  -0.21%  "

# Row 23
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").Value = "0.9992"
$ws.Range("E23").Value = "  -0.12%  "

# Row 24
$ws.Range("D24").Value = "6.082"
$ws.Range("E24").Value = "  -0.42%  "

# Row 25
$ws.Range("D25").Value = "9.224"
$ws.Range("E25").Value = "  -2.44%  "

# Row 26
$ws.Range("D26").Value = "163.32"
$ws.Range("E26").Value = "  +0.05%  "

# Row 27
$ws.Range("D27").Value = "18.51"
$ws.Range("E27").Value = "  +1.51%  "

# Row 28
$ws.Range("D28").Value = "1.917"
$ws.Range("E28").Value = "  -0.47%  "

# Row 29
$ws.Range("E29").Value = "  -2.81%  "

# Row 30
$ws.Range("D30").Value = "0.09194"
$ws.Range("E30").Value = "  +0.56%  "

# Row 31
$ws.Range("D31").Value = "4.240"
$ws.Range("E31").Value = "  -2.60%  "

# Row 32
$ws.Range("D32").Value = "3.956"
$ws.Range("E32").Value = "  -1.89%  "

# Row 33
$ws.Range("D33").Value = "0.05013"
$ws.Range("E33").Value = "  -3.51%  "

# Row 34
$ws.Range("D34").Value = "0.7392"
$ws.Range("E34").Value = "  -0.02%  "

# Row 35
$ws.Range("D35").Value = "1.142"
$ws.Range("E35").Value = "  +4.12%  "

# Row 36
$ws.Range("D36").Value = "2.686"
$ws.Range("E36").Value = "  -1.21%  "

# Row 37
$ws.Range("E37").Value = "  -0.05%  "

# Row 38
$ws.Range("D38").Value = "2.610"
$ws.Range("E38").Value = "  -2.02%  "

# Row 39
$ws.Range("D39").Value = "0.9064"
$ws.Range("E39").Value = "  -1.59%  "

# Row 40
$ws.Range("D40").Value = "2.054"
$ws.Range("E40").Value = "  +0.96%  "

# Row 41
$ws.Range("D41").Value = "5.958"
$ws.Range("E41").Value = "  +0.44%  "

# Row 42
$ws.Range("D42").Value = "106.51"
$ws.Range("E42").Value = "  +0.35%  "

# Row 43
$ws.Range("D43").Value = "0.4257"
$ws.Range("E43").Value = "  -3.03%  "

# Row 44
$ws.Range("D44").Value = "0.9985"
$ws.Range("E44").Value = "  +0.68%  "

# Row 45
$ws.Range("D45").Value = "7.392"
$ws.Range("E45").Value = "  -2.12%  "

# Row 46
$ws.Range("D46").Value = "0.1319"
$ws.Range("E46").Value = "  -3.61%  "

# Row 47
$ws.Range("D47").Value = "1.558"
$ws.Range("E47").Value = "  +11.39%  "

# Row 48
$ws.Range("D48").Value = "64.12"
$ws.Range("E48").Value = "  -6.12%  "

# Row 49
$ws.Range("D49").Value = "34.19"
$ws.Range("E49").Value = "  -2.07%  "

# Row 50
$ws.Range("D50").Value = "8.714"
$ws.Range("E50").Value = "  -2.49%  "

# Row 51
$ws.Range("D51").Value = "0.05656"
$ws.Range("E51").Value = "  -2.84%  "

# Restore the original (default) cell formatting/style now that the text
# values are locked in, so no stray number-format styling is left behind.
$editRange.ClearFormats()
